$wb = $excel.ActiveWorkbook

function Set-TextCell($cell, [string]$text) {
    # Force the cell to be written as a text value (t="inlineStr"/shared string)
    # rather than being auto-coerced into a number, while keeping the cell's
    # style back at the default (no explicit style) once done.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 1) "总计" sheet: insert a new "2022-Q3" row, pushing the existing "2022-Q1"
#    and "2021-Q3" rows down by one.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Create row 4 with the same formatting as row 3 (which already carries the
# bold/boxed "index" style used in column A), then fill in the former row 3
# values ("2021-Q3").
$summary.Range("A3").Copy()
$summary.Range("A4").PasteSpecial(-4122)  # xlPasteFormats
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(4, 2).Value = "2021-Q3"
$summary.Cells.Item(4, 3).Value = 2
$summary.Cells.Item(4, 4).Value = 0.13

# Row 3 becomes the old row 2 ("2022-Q1") values.
$summary.Cells.Item(3, 2).Value = "2022-Q1"
$summary.Cells.Item(3, 3).Value = 1
$summary.Cells.Item(3, 4).Value = 0.01

# Row 2 becomes the brand new "2022-Q3" values.
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 1
$summary.Cells.Item(2, 4).Value = 0.17

# ---------------------------------------------------------------------------
# 2) Add a new "2022-Q3" detail sheet (as a copy of "2022-Q1" so it keeps the
#    exact same layout/styling), placed right before "2022-Q1", and fill it
#    in with the new fund data.
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Copy($q1)
$q3 = $wb.Worksheets.Item("2022-Q1 (2)")
$q3.Name = "2022-Q3"

Set-TextCell $q3.Cells.Item(2, 2) "090013"
Set-TextCell $q3.Cells.Item(2, 3) "大成竞争优势混合"
Set-TextCell $q3.Cells.Item(2, 4) "6.88"
Set-TextCell $q3.Cells.Item(2, 5) "61.00"
Set-TextCell $q3.Cells.Item(2, 6) "2.54"
Set-TextCell $q3.Cells.Item(2, 7) "0.1748"
$q3.Cells.Item(2, 8).Value = 9

# Restore the originally-active tab ("2021-Q3") since copying a sheet makes
# the new copy active.
$wb.Worksheets.Item("2021-Q3").Activate()
